$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.198148666666667
$ws.Range("H2").Value = 15.594446
$ws.Range("I2").Value = 0.06659084221819957
$ws.Range("J2").Value = 0.06659084221819957
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06778666666666666
$ws.Range("N2").Value = 0.20336
$ws.Range("O2").Value = 0.01026668284214455
$ws.Range("P2").Value = 0.01026668284214455
$ws.Range("Q2").Value = 0.3523651709511111
$ws.Range("R2").Value = 3.17128653856
$ws.Range("S2").Value = 0.0006836670572455445
$ws.Range("T2").Value = 0.0006836670572455445

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.198148666666667
$ws.Range("H3").Value = 15.594446
$ws.Range("I3").Value = 0.06659084221819957
$ws.Range("J3").Value = 0.06659084221819957
$ws.Range("O3").Value = 0.01567037284022157
$ws.Range("P3").Value = 0.01567037284022157
$ws.Range("Q3").Value = 0.5378264517966667
$ws.Range("R3").Value = 4.84043806617
$ws.Range("S3").Value = 0.001043503325303554
$ws.Range("T3").Value = 0.001043503325303554

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.198148666666667
$ws.Range("H4").Value = 15.594446
$ws.Range("I4").Value = 0.06659084221819957
$ws.Range("J4").Value = 0.06659084221819957
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08443199999999999
$ws.Range("N4").Value = 0.253296
$ws.Range("O4").Value = 0.01278771487600239
$ws.Range("P4").Value = 0.01278771487600239
$ws.Range("Q4").Value = 0.4388900882240001
$ws.Range("R4").Value = 3.950010794016
$ws.Range("S4").Value = 0.0008515447036391986
$ws.Range("T4").Value = 0.0008515447036391985

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.198148666666667
$ws.Range("H5").Value = 15.594446
$ws.Range("I5").Value = 0.06659084221819957
$ws.Range("J5").Value = 0.06659084221819957
$ws.Range("M5").Value = 6.346903333333334
$ws.Range("N5").Value = 19.04071
$ws.Range("O5").Value = 0.9612752294416316
$ws.Range("P5").Value = 0.9612752294416316
$ws.Range("Q5").Value = 32.99214709962889
$ws.Range("R5").Value = 296.92932389666
$ws.Range("S5").Value = 0.06401212713201128
$ws.Range("T5").Value = 0.06401212713201128

# Row 6
$ws.Range("I6").Value = 0.1664471631553919
$ws.Range("J6").Value = 0.1664471631553919
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.06778666666666666
$ws.Range("N6").Value = 0.20336
$ws.Range("O6").Value = 0.01026668284214455
$ws.Range("P6").Value = 0.01026668284214455
$ws.Range("Q6").Value = 0.8807544873422222
$ws.Range("R6").Value = 7.92679038608
$ws.Range("S6").Value = 0.001708860234091097
$ws.Range("T6").Value = 0.001708860234091097

# Row 7
$ws.Range("I7").Value = 0.1664471631553919
$ws.Range("J7").Value = 0.1664471631553919
$ws.Range("O7").Value = 0.01567037284022157
$ws.Range("P7").Value = 0.01567037284022157
$ws.Range("S7").Value = 0.002608289104842182
$ws.Range("T7").Value = 0.002608289104842182

# Row 8
$ws.Range("I8").Value = 0.1664471631553919
$ws.Range("J8").Value = 0.1664471631553919
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08443199999999999
$ws.Range("N8").Value = 0.253296
$ws.Range("O8").Value = 0.01278771487600239
$ws.Range("P8").Value = 0.01278771487600239
$ws.Range("Q8").Value = 1.097027874832
$ws.Range("R8").Value = 9.873250873487999
$ws.Range("S8").Value = 0.002128478864350602
$ws.Range("T8").Value = 0.002128478864350602

# Row 9
$ws.Range("I9").Value = 0.1664471631553919
$ws.Range("J9").Value = 0.1664471631553919
$ws.Range("M9").Value = 6.346903333333334
$ws.Range("N9").Value = 19.04071
$ws.Range("O9").Value = 0.9612752294416316
$ws.Range("P9").Value = 0.9612752294416316
$ws.Range("Q9").Value = 82.46553292034778
$ws.Range("R9").Value = 742.18979628313
$ws.Range("S9").Value = 0.1600015349521081
$ws.Range("T9").Value = 0.1600015349521081

# Row 10
$ws.Range("G10").Value = 59.72155033333333
$ws.Range("H10").Value = 179.164651
$ws.Range("I10").Value = 0.7650624463235045
$ws.Range("J10").Value = 0.7650624463235045
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.06778666666666666
$ws.Range("N10").Value = 0.20336
$ws.Range("O10").Value = 0.01026668284214455
$ws.Range("P10").Value = 0.01026668284214455
$ws.Range("Q10").Value = 4.048324825262222
$ws.Range("R10").Value = 36.43492342736
$ws.Range("S10").Value = 0.007854653490838661
$ws.Range("T10").Value = 0.007854653490838661

# Row 11
$ws.Range("G11").Value = 59.72155033333333
$ws.Range("H11").Value = 179.164651
$ws.Range("I11").Value = 0.7650624463235045
$ws.Range("J11").Value = 0.7650624463235045
$ws.Range("O11").Value = 0.01567037284022157
$ws.Range("P11").Value = 0.01567037284022157
$ws.Range("Q11").Value = 6.179090205238333
$ws.Range("R11").Value = 55.61181184714499
$ws.Range("S11").Value = 0.01198881377994132
$ws.Range("T11").Value = 0.01198881377994132

# Row 12
$ws.Range("G12").Value = 59.72155033333333
$ws.Range("H12").Value = 179.164651
$ws.Range("I12").Value = 0.7650624463235045
$ws.Range("J12").Value = 0.7650624463235045
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.08443199999999999
$ws.Range("N12").Value = 0.253296
$ws.Range("O12").Value = 0.01278771487600239
$ws.Range("P12").Value = 0.01278771487600239
$ws.Range("Q12").Value = 5.042409937744
$ws.Range("R12").Value = 45.38168943969599
$ws.Range("S12").Value = 0.00978340042592186
$ws.Range("T12").Value = 0.009783400425921858

# Row 13
$ws.Range("G13").Value = 59.72155033333333
$ws.Range("H13").Value = 179.164651
$ws.Range("I13").Value = 0.7650624463235045
$ws.Range("J13").Value = 0.7650624463235045
$ws.Range("M13").Value = 6.346903333333334
$ws.Range("N13").Value = 19.04071
$ws.Range("O13").Value = 0.9612752294416316
$ws.Range("P13").Value = 0.9612752294416316
$ws.Range("Q13").Value = 379.0469068824678
$ws.Range("R13").Value = 3411.42216194221
$ws.Range("S13").Value = 0.7354355786268028
$ws.Range("T13").Value = 0.7354355786268028

# Row 14
$ws.Range("G14").Value = 0.1482806666666667
$ws.Range("H14").Value = 0.444842
$ws.Range("I14").Value = 0.001899548302904017
$ws.Range("J14").Value = 0.001899548302904017
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.06778666666666666
$ws.Range("N14").Value = 0.20336
$ws.Range("O14").Value = 0.01026668284214455
$ws.Range("P14").Value = 0.01026668284214455
$ws.Range("Q14").Value = 0.01005145212444444
$ws.Range("R14").Value = 0.09046306912
$ws.Range("S14").Value = 0.00001950205996924947
$ws.Range("T14").Value = 0.00001950205996924947

# Row 15
$ws.Range("G15").Value = 0.1482806666666667
$ws.Range("H15").Value = 0.444842
$ws.Range("I15").Value = 0.001899548302904017
$ws.Range("J15").Value = 0.001899548302904017
$ws.Range("O15").Value = 0.01567037284022157
$ws.Range("P15").Value = 0.01567037284022157
$ws.Range("Q15").Value = 0.01534185917666667
$ws.Range("R15").Value = 0.13807673259
$ws.Range("S15").Value = 0.00002976663013451607
$ws.Range("T15").Value = 0.00002976663013451607

# Row 16
$ws.Range("G16").Value = 0.1482806666666667
$ws.Range("H16").Value = 0.444842
$ws.Range("I16").Value = 0.001899548302904017
$ws.Range("J16").Value = 0.001899548302904017
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.08443199999999999
$ws.Range("N16").Value = 0.253296
$ws.Range("O16").Value = 0.01278771487600239
$ws.Range("P16").Value = 0.01278771487600239
$ws.Range("Q16").Value = 0.012519633248
$ws.Range("R16").Value = 0.112676699232
$ws.Range("S16").Value = 0.00002429088209073079
$ws.Range("T16").Value = 0.00002429088209073078

# Row 17
$ws.Range("G17").Value = 0.1482806666666667
$ws.Range("H17").Value = 0.444842
$ws.Range("I17").Value = 0.001899548302904017
$ws.Range("J17").Value = 0.001899548302904017
$ws.Range("M17").Value = 6.346903333333334
$ws.Range("N17").Value = 19.04071
$ws.Range("O17").Value = 0.9612752294416316
$ws.Range("P17").Value = 0.9612752294416316
$ws.Range("Q17").Value = 0.9411230575355556
$ws.Range("R17").Value = 8.470107517820001
$ws.Range("S17").Value = 0.00182598873070952
$ws.Range("T17").Value = 0.00182598873070952
